$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Fix the wording in the trimmer note: "...lower and top left pads..." -> "...lower and top right pads..."
$ws.Range("B16").Value = "* Trimmers TM1-4 will let you easily change LED colour and adjust the brightness. If you feel like calculating the current limiting`nresistance according to your LED specifications you may solder a fixed resistor connecting the lower and top right pads`nof the three smd pads instead. Then leave the third pad unmounted."

# Scroll the view down so row 8 is the top visible row, then select E10
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E10").Select()
